$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.150445461273193
$ws.Range("B1").Value = 3.084508419036865
$ws.Range("C1").Value = 2.730147838592529
$ws.Range("D1").Value = 1.682782888412476
$ws.Range("E1").Value = 0.8727366328239441
